$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.035.78'
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("D3").Value = '1.597.73'
$ws.Range("E3").Value = '  -0.08%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '302.17'
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3773'
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3631'
$ws.Range("E8").Value = '  -1.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '50.84'
$ws.Range("E9").Value = '  +3.92%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.248'
$ws.Range("E10").Value = '  -2.46%  '
$ws.Range("E11").Value = '  -0.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08132'
$ws.Range("E12").Value = '  +0.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.31'
$ws.Range("E13").Value = '  -2.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.561'
$ws.Range("E14").Value = '  -1.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.355'
$ws.Range("E15").Value = '  -2.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001240'
$ws.Range("E16").Value = '  -2.00%  '
$ws.Range("D17").Value = '1.597.41'
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.26'
$ws.Range("E18").Value = '  +0.38%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06843'
$ws.Range("E19").Value = '  +0.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.13'
$ws.Range("E20").Value = '  -2.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.494'
$ws.Range("E21").Value = '  -1.93%  '
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.97'
$ws.Range("E23").Value = '  -1.46%  '
$ws.Range("D24").Value = '23.028.11'
$ws.Range("E24").Value = '  -0.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.375'
$ws.Range("E25").Value = '  +0.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.789'
$ws.Range("E26").Value = '  -6.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.08'
$ws.Range("E27").Value = '  -0.39%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '149.15'
$ws.Range("E28").Value = '  -1.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.263'
$ws.Range("E29").Value = '  -0.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '135.00'
$ws.Range("E30").Value = '  +1.80%  '
$ws.Range("E31").Value = '  -4.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.736'
$ws.Range("E32").Value = '  -5.70%  '
$ws.Range("D33").Value = '1.773.81'
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9563'
$ws.Range("E34").Value = '  -1.89%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07514'
$ws.Range("E35").Value = '  -3.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02706'
$ws.Range("E36").Value = '  -2.96%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.185'
$ws.Range("E37").Value = '  -1.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.13'
$ws.Range("E38").Value = '  -0.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2512'
$ws.Range("E39").Value = '  -1.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08812'
$ws.Range("E40").Value = '  -0.79%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.357'
$ws.Range("E41").Value = '  -1.77%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7018'
$ws.Range("E42").Value = '  -2.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.36'
$ws.Range("E43").Value = '  -3.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.29'
$ws.Range("E44").Value = '  -5.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6564'
$ws.Range("E45").Value = '  -1.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.998'
$ws.Range("E46").Value = '  +0.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.268'
$ws.Range("E47").Value = '  -2.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '131.80'
$ws.Range("E48").Value = '  +0.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07928'
$ws.Range("E49").Value = '  -1.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.216'
$ws.Range("E50").Value = '  +3.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.225'
$ws.Range("E51").Value = '  +2.81%  '
